$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-01-27 Monday"; new = "2025-01-28 Tuesday"},
    @{old = "29×55=1595"; new = "73×24=1752"},
    @{old = "17×69=1173"; new = "67×21=1407"},
    @{old = "16×91=1456"; new = "84×21=1764"},
    @{old = "63×23=1449"; new = "20×70=1400"},
    @{old = "31×92=2852"; new = "91×33=3003"},
    @{old = "57×30=1710"; new = "94×13=1222"},
    @{old = "72×36=2592"; new = "83×14=1162"},
    @{old = "46×84=3864"; new = "44×53=2332"},
    @{old = "24×84=2016"; new = "82×82=6724"},
    @{old = "37×49=1813"; new = "70×83=5810"},
    @{old = "99×21=2079"; new = "92×68=6256"},
    @{old = "25×58=1450"; new = "61×15=915"},
    @{old = "46×28=1288"; new = "11×64=704"},
    @{old = "21×17=357"; new = "37×65=2405"},
    @{old = "76×39=2964"; new = "52×53=2756"},
    @{old = "99×83=8217"; new = "40×51=2040"},
    @{old = "65×89=5785"; new = "98×90=8820"},
    @{old = "50×79=3950"; new = "85×15=1275"},
    @{old = "89×92=8188"; new = "87×11=957"},
    @{old = "87×17=1479"; new = "91×42=3822"},
    @{old = "77×32=2464"; new = "96×40=3840"},
    @{old = "80×98=7840"; new = "65×60=3900"},
    @{old = "29×93=2697"; new = "15×14=210"},
    @{old = "85×47=3995"; new = "54×80=4320"},
    @{old = "53×64=3392"; new = "95×68=6460"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
